$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 1).Value = "id_DK_nan_PV"
$ws.Cells.Item(3, 2).Value = "PV_DK"
$ws.Cells.Item(4, 1).Value = "id_DK_nan_ROR"
$ws.Cells.Item(4, 2).Value = "ROR_DK"
$ws.Cells.Item(5, 1).Value = "id_DK_nan_WL"
$ws.Cells.Item(5, 2).Value = "WL_DK"
$ws.Cells.Item(6, 1).Value = "id_DK_nan_WS"
$ws.Cells.Item(6, 2).Value = "WS_DK"
$ws.Cells.Item(7, 1).Value = "id_DK_Central_BH_Biogas"
$ws.Cells.Item(7, 2).Value = "Standard"
$ws.Cells.Item(8, 1).Value = "id_DK_Central_BP_Biogas"
$ws.Cells.Item(8, 2).Value = "Standard"
$ws.Cells.Item(9, 1).Value = "id_DK_Central_IndustryH_Biogas"
$ws.Cells.Item(9, 2).Value = "Standard"
$ws.Cells.Item(10, 1).Value = "id_DK_Central_BH_Biomass"
$ws.Cells.Item(10, 2).Value = "Standard"
$ws.Cells.Item(11, 1).Value = "id_DK_Central_BP_Biomass"
$ws.Cells.Item(11, 2).Value = "Standard"
$ws.Cells.Item(12, 1).Value = "id_DK_Central_IndustryH_Biomass"
$ws.Cells.Item(12, 2).Value = "Standard"
$ws.Cells.Item(13, 1).Value = "id_DK_Central_BP_Coal"
$ws.Cells.Item(13, 2).Value = "Standard"
$ws.Cells.Item(14, 1).Value = "id_DK_Central_BH_Natgas"
$ws.Cells.Item(14, 2).Value = "Standard"
$ws.Cells.Item(15, 1).Value = "id_DK_Central_BP_Natgas"
$ws.Cells.Item(15, 2).Value = "Standard"
$ws.Cells.Item(16, 1).Value = "id_DK_Central_IndustryH_Natgas"
$ws.Cells.Item(16, 2).Value = "Standard"
$ws.Cells.Item(17, 1).Value = "id_DK_Central_BH_Oil"
$ws.Cells.Item(17, 2).Value = "Standard"
$ws.Cells.Item(18, 1).Value = "id_DK_Central_BP_Oil"
$ws.Cells.Item(18, 2).Value = "Standard"
$ws.Cells.Item(19, 1).Value = "id_DK_Central_IndustryH_Oil"
$ws.Cells.Item(19, 2).Value = "Standard"
$ws.Cells.Item(20, 1).Value = "id_DK_Central_BH_Waste"
$ws.Cells.Item(20, 2).Value = "Standard"
$ws.Cells.Item(21, 1).Value = "id_DK_Central_BP_Waste"
$ws.Cells.Item(21, 2).Value = "Standard"
$ws.Cells.Item(22, 1).Value = "id_DK_Central_EP"
$ws.Cells.Item(22, 2).Value = "Standard"
$ws.Cells.Item(23, 1).Value = "id_DK_Central_GT"
$ws.Cells.Item(23, 2).Value = "Standard"
$ws.Cells.Item(24, 1).Value = "id_DK_Central_HPstandard"
$ws.Cells.Item(24, 2).Value = "Standard"
$ws.Cells.Item(25, 1).Value = "id_DK_Central_HPsurplusheat"
$ws.Cells.Item(25, 2).Value = "Standard"
$ws.Cells.Item(26, 1).Value = "id_DK_Central_IH"
$ws.Cells.Item(26, 2).Value = "Standard"
$ws.Cells.Item(27, 1).Value = "id_DK_Central_IndustryH"
$ws.Cells.Item(27, 2).Value = "Standard"
$ws.Cells.Item(28, 1).Value = "id_DK_nan_CD_Biogas"
$ws.Cells.Item(28, 2).Value = "Standard"
$ws.Cells.Item(29, 1).Value = "id_DK_nan_IndustryE_Biogas"
$ws.Cells.Item(29, 2).Value = "Standard"
$ws.Cells.Item(30, 1).Value = "id_DK_nan_IndustryE_Biomass"
$ws.Cells.Item(30, 2).Value = "Standard"
$ws.Cells.Item(31, 1).Value = "id_DK_nan_CD_Coal"
$ws.Cells.Item(31, 2).Value = "Standard"
$ws.Cells.Item(32, 1).Value = "id_DK_nan_CD_Natgas"
$ws.Cells.Item(32, 2).Value = "Standard"
$ws.Cells.Item(33, 1).Value = "id_DK_nan_IndustryE_Natgas"
$ws.Cells.Item(33, 2).Value = "Standard"
$ws.Cells.Item(34, 1).Value = "id_DK_nan_CD_Oil"
$ws.Cells.Item(34, 2).Value = "Standard"
$ws.Cells.Item(35, 1).Value = "id_DK_nan_IndustryE_Oil"
$ws.Cells.Item(35, 2).Value = "Standard"

# Remove now-obsolete trailing rows (table shrank from 54 to 35 rows)
$ws.Range("A36:B54").Clear()

